# Add 2022-Q3 data:
#  - Insert a new "2022-Q3" row into the "总计" (summary) sheet, shifting
#    the existing quarter rows down by one.
#  - Insert a new "2022-Q3" worksheet (before "2022-Q2") holding the
#    per-fund detail rows for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new row 2 for 2022-Q3.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Bring row 2's "A" (index) cell formatting into line with the other
# index cells (style carries border/alignment), then strip the
# formatting that Insert() auto-propagated into B2:D2 from row 1.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0

# Renumber the index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet with the fund detail rows, placed right
#    before the existing "2022-Q2" sheet. Cloning "2022-Q2" keeps the
#    header row / styles / sheetPr / page margins identical, then we
#    drop its data rows beyond row 3 and overwrite rows 2-3.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2Index = $q2.Index
$q2.Copy($q2)
# The copy is inserted at $q2's old position; the original "2022-Q2" is
# pushed one slot later. Grab the copy by that original index.
$q3 = $wb.Worksheets.Item($q2Index)
$q3.Name = "2022-Q3"

# Drop the old 2022-Q2 data rows beyond the 2 new detail rows.
$q3.Range("4:12").Delete()

# Row 2 — 瑞达策略优选混合A
$q3.Range("A2").Value = 0
$q3.Range("B2:G2").NumberFormat = "@"
$q3.Range("B2").Value = "015694"
$q3.Range("C2").Value = "瑞达策略优选混合A"
$q3.Range("D2").Value = "0.09"
$q3.Range("E2").Value = "67.87"
$q3.Range("F2").Value = "2.79"
$q3.Range("G2").Value = "0.0025"
$q3.Range("B2:G2").ClearFormats()
$q3.Range("H2").Value = 5

# Row 3 — 瑞达策略优选混合C
$q3.Range("A3").Value = 1
$q3.Range("B3:F3").NumberFormat = "@"
$q3.Range("B3").Value = "015695"
$q3.Range("C3").Value = "瑞达策略优选混合C"
$q3.Range("D3").Value = "0.00"
$q3.Range("E3").Value = "67.87"
$q3.Range("F3").Value = "2.79"
$q3.Range("B3:F3").ClearFormats()
$q3.Range("G3").Value = 0
$q3.Range("H3").Value = 5

# Copying a sheet makes it the active one; restore the original active
# tab ("2021-Q1", the last sheet) so the view state is unchanged.
$wb.Worksheets.Item("2021-Q1").Activate()
